# Fix formatting issues introduced when scraping floating point numbers.
# 1) A handful of "Razon social" entries used a comma to separate multiple
#    co-contracted people/parties; those commas are changed to periods (and
#    a couple of "S.H." abbreviations are normalised to "SH") so they are no
#    longer confused with the decimal/thousands separators used elsewhere.
# 2) The "Importe" column (H) held amounts as text using the es-AR style
#    ("1.234.567,89" -- "." thousands separator, "," decimal separator).
#    They are rewritten as plain "1234567.89" style numbers-as-text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Razon social fixes -------------------------------------------------
$ws.Range("E95").Value  = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E178").Value = "EDICIONES NATIVA SH DE ESCOBAR JORGE. MARTINEZ ALFREDO. PIZIGHINI CARLOS L Y R"
$ws.Range("E214").Value = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
$ws.Range("E227").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"

# --- 2) Importe (column H) reformatting ------------------------------------
# Values are stored as text, so keep the column formatted as Text before
# writing the new strings -- otherwise Excel would "helpfully" reinterpret
# them as numbers (and drop the trailing zeros we want to keep).
$ws.Range("H2:H303").NumberFormat = "@"

$importes = @(
  "155250.00",
  "17000.00",
  "93010.00",
  "4032000.00",
  "1374000.00",
  "1374000.00",
  "250000.00",
  "5908320.00",
  "3184800.00",
  "600.00",
  "19826.00",
  "365.00",
  "1540.00",
  "14810.00",
  "1504477.40",
  "587000.00",
  "315000.00",
  "390.00",
  "4200.00",
  "39194.00",
  "953809.04",
  "44206.00",
  "3100.00",
  "21400.00",
  "94590.00",
  "122880.00",
  "562833.08",
  "720.00",
  "3096.00",
  "4782.10",
  "92265.88",
  "24923.00",
  "48000.00",
  "117785.43",
  "9720.00",
  "39900.00",
  "15972.00",
  "3000.00",
  "13418.90",
  "6800.00",
  "4500.00",
  "2190.00",
  "67.49",
  "359600.00",
  "909.20",
  "247.94",
  "106000.00",
  "1657.82",
  "2891056.53",
  "151512.73",
  "1100.00",
  "23540.00",
  "24268.00",
  "433.25",
  "458429.71",
  "117240.00",
  "1050.00",
  "15792.75",
  "10045.02",
  "120.00",
  "4450.00",
  "15196.61",
  "33549.50",
  "102373.40",
  "4154.00",
  "16003.40",
  "90210.00",
  "1880.00",
  "1279551.69",
  "2490.00",
  "5523.00",
  "735.50",
  "12530.29",
  "4583.86",
  "105956.65",
  "33960.00",
  "53200.00",
  "55900.00",
  "30835.00",
  "1650.00",
  "3139.46",
  "8680.44",
  "5897.02",
  "15030.00",
  "144000.00",
  "1170.00",
  "6900.00",
  "172702.90",
  "850.00",
  "32900.00",
  "5700.00",
  "6600.00",
  "134940.00",
  "3360.00",
  "7195.00",
  "350.00",
  "543300.00",
  "19276.00",
  "2142.00",
  "6352.00",
  "18.75",
  "36700.48",
  "9.25",
  "5.20",
  "43288.80",
  "28748.12",
  "7500.00",
  "944.80",
  "6957.00",
  "30410.00",
  "125.00",
  "43904.40",
  "17360.00",
  "103256.75",
  "120.00",
  "2518500.00",
  "7340.00",
  "375.00",
  "2480.00",
  "2500.00",
  "11770.00",
  "1900.00",
  "1904.00",
  "6870.00",
  "8603.00",
  "26511.78",
  "2500.00",
  "193.91",
  "852.00",
  "1050.00",
  "13481.00",
  "477.69",
  "4072.67",
  "3289.00",
  "286.06",
  "5400.00",
  "615767.02",
  "32000.00",
  "338500.00",
  "62000.00",
  "61828.20",
  "89020.00",
  "6852.00",
  "12290.00",
  "7000.00",
  "17000.00",
  "10000.00",
  "6000.00",
  "42000.00",
  "129000.00",
  "16200.00",
  "13200.00",
  "194000.00",
  "444900.00",
  "46000.00",
  "15600.00",
  "7600.00",
  "49000.00",
  "19076.22",
  "8800.00",
  "23327.70",
  "13120.00",
  "51300.00",
  "7851.91",
  "14445.00",
  "13038000.00",
  "35805.00",
  "18000.00",
  "16500.00",
  "22000.00",
  "22000.00",
  "8000.00",
  "7000.00",
  "12000.00",
  "60000.00",
  "8500.00",
  "23728.77",
  "12000.00",
  "10000.00",
  "10000.00",
  "9000.00",
  "9000.00",
  "4000.00",
  "10000.00",
  "10000.00",
  "5000.00",
  "21000.00",
  "10000.00",
  "10000.00",
  "5000.00",
  "59901.50",
  "16000.00",
  "10000.00",
  "12000.00",
  "4500.00",
  "12000.00",
  "3500.00",
  "28000.00",
  "10000.00",
  "10000.00",
  "85498.88",
  "3500.00",
  "93320.00",
  "10000.00",
  "11000.00",
  "11000.00",
  "2420.00",
  "16424.00",
  "1200.00",
  "46250.00",
  "7200.00",
  "98.76",
  "7840.00",
  "2783.00",
  "137.98",
  "23180.00",
  "8700.00",
  "4168.24",
  "10580.00",
  "4561.05",
  "7260.00",
  "14850.00",
  "15050.00",
  "65651.00",
  "16806.79",
  "23940.00",
  "1117.65",
  "16740.90",
  "2800.00",
  "6017.43",
  "450.00",
  "9450.00",
  "30485.14",
  "24282.88",
  "15854.54",
  "80000.00",
  "40000.00",
  "40000.00",
  "40000.00",
  "80000.00",
  "40000.00",
  "55000.00",
  "40000.00",
  "40000.00",
  "80000.00",
  "80000.00",
  "71700.00",
  "222786.17",
  "6500.00",
  "13000.00",
  "239718.63",
  "4650.00",
  "199800.00",
  "20691.00",
  "815000.00",
  "258000.00",
  "4595095.62",
  "281850.00",
  "258000.00",
  "270400.00",
  "258000.00",
  "272400.00",
  "486000.00",
  "258000.00",
  "627450.00",
  "577000.00",
  "342600.00",
  "258000.00",
  "508000.00",
  "516000.00",
  "415000.00",
  "512400.00",
  "737400.00",
  "486000.00",
  "742300.00",
  "516000.00",
  "265970.00",
  "69700.00",
  "243932.22",
  "2422.33",
  "2193933.07",
  "9900.00",
  "1200.00",
  "186000.00",
  "4800.00",
  "15000.00",
  "11700.00",
  "10000.00",
  "36000.00",
  "13500.00",
  "56300.00",
  "597000.00",
  "1800.00",
  "116900.00",
  "17000.00",
  "6560.02",
  "5400.00",
  "300.00",
  "12000.00",
  "9700.00",
  "70225.00",
  "4750.00"
)

for ($i = 0; $i -lt $importes.Length; $i++) {
    $ws.Cells.Item($i + 2, 8).Value = $importes[$i]
}
